$d = $word.ActiveDocument

# Locate the paragraph that follows the "User Agent sempre ..." paragraph -
# that's the (currently empty) paragraph whose pPr gains <w:spacing w:after="0"/>
# and after which the new "Tabelas HTML5" section is inserted.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*User Agent sempre*exibido*conte*do*") {
        $anchor = $p
    }
}
$target = $anchor.Next()

# Hunk 1: add <w:spacing w:after="0"/> to that paragraph's pPr.
$target.SpaceAfter = 0

# Hunk 2: insert the new paragraphs right after it (before the final blank
# paragraph), built as literal OOXML so every run/paragraph-mark property
# (b, bCs, sz/szCs, spacing) matches exactly.
$newParagraphsXml = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Tabelas HTML5</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Não pode usar</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tabela no html5</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> para cria um site, usar tabela para poder criar uma tabela.</w:t></w:r></w:p>'

$package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Insert *before* the final (last) paragraph of the document, right after
# $target, by collapsing the range of the paragraph that follows $target
# to its start and inserting there.
$insertionPoint = $target.Next()
$insertRange = $insertionPoint.Range
$null = $insertRange.Collapse(1)
$null = $insertRange.InsertXML($package)

Write-Output "Paragraphs.Count=$($d.Paragraphs.Count)"
